$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename place cell A19 (Swamp -> In a swamp) first so new shared strings
# are appended in the same order as the target workbook.
$ws.Range("A19").Value = "In a swamp"

# Rename header cell A1 ("Name of the place" -> "Name of place when there")
$ws.Range("A1").Value = "Name of place when there"

# Column A's "best fit" width shrinks now that the longest entries changed
# ("Name of the place"/"Swamp" no longer drive the width).
$ws.Columns.Item(1).ColumnWidth = 28

# Update selection to match new active cell
$ws.Range("D20").Select() | Out-Null
